$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 62; existing rows 62-74 shift down to 63-75.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new record.
$ws.Cells.Item(62, 1).Value = 5
$ws.Cells.Item(62, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(62, 3).Value = 'Maule'
$ws.Cells.Item(62, 4).Value = 44505
$ws.Cells.Item(62, 5).Value = 7
$ws.Cells.Item(62, 6).Value = 100112013
$ws.Cells.Item(62, 7).Value = 'Alcachofa'
$ws.Cells.Item(62, 8).Value = 'Madrigal'
$ws.Cells.Item(62, 9).Value = 'Primera'
$ws.Cells.Item(62, 10).Value = 200
$ws.Cells.Item(62, 11).Value = 10000
$ws.Cells.Item(62, 12).Value = 10000
$ws.Cells.Item(62, 13).Value = 10000
$ws.Cells.Item(62, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(62, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(62, 16).Value = 250
$ws.Cells.Item(62, 17).Value = 40
$ws.Cells.Item(62, 18).Value = 'Hortaliza'
